$d = $word.ActiveDocument

# --- Title ---
$d.Content.Find.Execute("Unveiling the Enigmatic Depths of Space", $true, $false, $false, $false, $false, $true, 1, $false, "Unveiling the Marvels of Chemistry: A Journey Through the Molecular Realm", 1) | Out-Null

# --- Author name ---
$d.Content.Find.Execute(" Clarissa McPherson", $true, $false, $false, $false, $false, $true, 1, $false, " Eleanor Wilson", 1) | Out-Null

# --- Email line (paragraph 3): replace the three runs precisely ---
$emailPara = $d.Paragraphs(3)
$emailStart = $emailPara.Range.Start
# "c" -> "eleanor" (single leading character run; positional to avoid matching letters elsewhere)
$d.Range($emailStart, $emailStart + 1).Text = "eleanor"
$emailPara = $d.Paragraphs(3)
$emailPara.Range.Find.Execute("mcpherson@spatiumobservatory", $true, $false, $false, $false, $false, $true, 0, $false, "wilson098@college", 1) | Out-Null
$emailPara = $d.Paragraphs(3)
$emailPara.Range.Find.Execute("org", $true, $false, $false, $false, $false, $true, 0, $false, "edu", 1) | Out-Null

# --- Body paragraph (Summary intro block) ---
$d.Content.Find.Execute("The cosmic canvas above holds an inexhaustible reservoir of mysteries, beckoning humanity to embark on an unceasing quest for knowledge and understanding", $true, $false, $false, $false, $false, $true, 1, $false, "Within the realm of chemistry, we embark on an awe-inspiring journey through the microscopic world of elements and molecules", 1) | Out-Null
$d.Content.Find.Execute(" From the earliest civilizations gazing up in awe at the celestial tapestry to the modern era of sophisticated observatories, we continue to unravel the enigmas that shroud the vast expanses of space", $true, $false, $false, $false, $false, $true, 1, $false, " Chemistry, the science that delves into the composition, structure, properties, and changes of matter, holds profound significance in our lives, shaping the world we inhabit", 1) | Out-Null
$d.Content.Find.Execute(" The pursuit of space exploration harnesses groundbreaking technologies, enabling us to traverse the cosmos, study extraterrestrial bodies, and search for life beyond Earth. As we delve deeper into the cosmic abyss, we witness celestial phenomena of immense beauty and complexity, challenging our comprehension of the universe's origins, evolution, and ultimate fate", $true, $false, $false, $false, $false, $true, 1, $false, " From the air we breathe to the food we consume, chemistry plays an integral role in understanding and manipulating the intricate processes that govern our natural and technological landscapes", 1) | Out-Null
$d.Content.Find.Execute("Space exploration serves as a catalyst for scientific advancements, propelling us to develop cutting-edge instruments, delve into exotic realms of physics, and refine our understanding of fundamental cosmic laws", $true, $false, $false, $false, $false, $true, 1, $false, "Unraveling the secrets of chemistry unveils a captivating tapestry of interactions between atoms and molecules", 1) | Out-Null
$d.Content.Find.Execute(" Missions to distant worlds reveal hidden treasures--diverse planetary systems, awe-inspiring moons, and potentially habitable exoplanets", $true, $false, $false, $false, $false, $true, 1, $false, " We explore the fundamental building blocks of matter, tracing the evolution of elements from their primordial origins in the cosmic crucible to their diverse manifestations in the world around us", 1) | Out-Null
$d.Content.Find.Execute(" Furthermore, the study of space unveils profound insights into the formation and evolution of stars, galaxies, and the cosmos as a whole. Exploring the cosmos is not merely an intellectual endeavor; it enriches our cultural tapestry, inspires artistic expression, and ignites our existential ponderings about our place in the vastness of the universe", $true, $false, $false, $false, $false, $true, 1, $false, " By unraveling the intricacies of chemical reactions, we decipher the mechanisms that drive biological processes, unlock the potential of materials, and devise innovative technologies that shape our modern existence", 1) | Out-Null
$d.Content.Find.Execute("With each successful mission and discovery, humanity gains not only knowledge but also a renewed sense of wonder and humility", $true, $false, $false, $false, $false, $true, 1, $false, "Chemistry, in its boundless applications, touches every facet of human endeavor", 1) | Out-Null
$d.Content.Find.Execute(" We realize that our planet is but a small speck in an infinite expanse, interconnected with countless celestial bodies", $true, $false, $false, $false, $false, $true, 1, $false, " It empowers us to harness the energy stored in fossil fuels and renewable sources, providing the lifeblood that drives our economies and powers our lives", 1) | Out-Null
$d.Content.Find.Execute(" The exploration of space fosters global collaboration, uniting scientists, engineers, and visionaries from diverse backgrounds in a common pursuit of understanding our place in the cosmos", $true, $false, $false, $false, $false, $true, 1, $false, " It enables the synthesis of pharmaceuticals that combat disease and alleviate suffering, safeguarding our health and well-being", 1) | Out-Null
$d.Content.Find.Execute(" As we embark on this extraordinary journey of discovery, we embrace the unknown, relentlessly pushing the boundaries of human knowledge and imagination", $true, $false, $false, $false, $false, $true, 1, $false, " Through the transformative power of chemistry, we engineer materials with extraordinary properties, paving the way for advancements in electronics, transportation, and construction, reshaping our built environment", 1) | Out-Null

# --- Summary paragraph ---
$d.Content.Find.Execute("Humankind's exploration of space is an ongoing saga of scientific discovery, technological innovation, and profound existential contemplation", $true, $false, $false, $false, $false, $true, 1, $false, "Chemistry, the study of the composition, structure, properties, and changes of matter, unveils the intricate world of elements and molecules that shape our lives", 1) | Out-Null
$d.Content.Find.Execute(" We venture into the cosmic frontier, unraveling the enigmas of the universe and gaining invaluable insights into our place within it", $true, $false, $false, $false, $false, $true, 1, $false, " It unravels the secrets of chemical reactions, empowering us to decipher biological processes, unlock the potential of materials, and engineer innovative technologies", 1) | Out-Null
$d.Content.Find.Execute(" Space exploration transcends national boundaries, fostering global collaboration and inspiring generations to reach for the stars. It is a testament to our insatiable curiosity and unyielding desire to comprehend the vastness and wonder that surrounds us", $true, $false, $false, $false, $false, $true, 1, $false, " From energy production to healthcare advancements, and the development of revolutionary materials, chemistry's boundless applications touch every aspect of human endeavor, driving progress and transforming our world", 1) | Out-Null

# --- Add a new empty paragraph at the very end of the body ---
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

Write-Output "edit complete"
